$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "byTrial"

# New sheet, inserted right after byTrial.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "byAcquisition"

# Header row (reuse the same shared strings / order as byTrial).
$headers = @("blink feature", "r", "mean slope", "t-stat", "p")
for ($c = 0; $c -lt 5; $c++) {
    $ws2.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Data rows: feature, r, mean slope, t-stat, p
$rows = @(
    @("auc",                0.4773,   -0.0095,  -2.5298,  0.024),
    @("latency",            -0.0232,   0.0188,   1.097,    0.2912),
    @("timeUnder20",        -0.1304,  -0.0252,  -4.4953,   0.0005),
    @("openTime",           -0.0568,  -0.0772,  -4.3287,   0.0007),
    @("initialVelocity",     0.47,     0.2374,   0.2043,   0.841),
    @("closeTime",           0.2789,   0.1166,   0.194,    0.849),
    @("maxClosingVelocity",  0.4304,  -0.3162,  -2.3398,   0.0346),
    @("maxOpeningVelocity",  0.0583,   7.6755,   1.6463,   0.1219),
    @("blinkRate",           0.7987,  -0.0267,  -2.3955,   0.0311)
)

$r = 2
foreach ($row in $rows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Copy header formatting (bold + centered) from byTrial's header row, reusing the existing style.
$ws1.Range("A1:E1").Copy() | Out-Null
$ws2.Range("A1:E1").PasteSpecial(-4122) | Out-Null

# Column widths, matching byTrial's.
$ws2.Columns.Item(1).ColumnWidth = $ws1.Columns.Item(1).ColumnWidth
$ws2.Columns.Item(2).ColumnWidth = $ws1.Columns.Item(2).ColumnWidth

# Selections / active tab: byAcquisition tab is selected, with D20 selected there;
# byTrial keeps B41 selected but is no longer the active tab.
$ws2.Range("D20").Select() | Out-Null
$ws1.Range("B41").Select() | Out-Null
$ws2.Activate() | Out-Null

Write-Host "done"
